$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "pdi" (column I) values with more precise measurements
$ws.Range("I7").Value = 0.289
$ws.Range("I8").Value = 0.161
$ws.Range("I14").Value = 0.129
$ws.Range("I16").Value = 0.151
$ws.Range("I18").Value = 0.181
$ws.Range("I24").Value = 0.179
$ws.Range("I30").Value = 0.196

# Update the view scroll position and active selection
$excel.ActiveWindow.ScrollRow = 6
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("I31").Select()
